$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-08-20 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-21 Monday", 2) | Out-Null

# Update the 20x5 answer table
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "86-48=38"
$t.Cell(1, 2).Range.Text = "26+38=64"
$t.Cell(1, 3).Range.Text = "63+34=97"
$t.Cell(1, 4).Range.Text = "96-30=66"
$t.Cell(1, 5).Range.Text = "59+23=82"
$t.Cell(2, 1).Range.Text = "72-3=69"
$t.Cell(2, 2).Range.Text = "44-7=37"
$t.Cell(2, 3).Range.Text = "47-11=36"
$t.Cell(2, 4).Range.Text = "35+57=92"
$t.Cell(2, 5).Range.Text = "27+64=91"
$t.Cell(3, 1).Range.Text = "71-43=28"
$t.Cell(3, 2).Range.Text = "93-76=17"
$t.Cell(3, 3).Range.Text = "86-81=5"
$t.Cell(3, 4).Range.Text = "74-37=37"
$t.Cell(3, 5).Range.Text = "99-94=5"
$t.Cell(4, 1).Range.Text = "78-68=10"
$t.Cell(4, 2).Range.Text = "16+83=99"
$t.Cell(4, 3).Range.Text = "66-15=51"
$t.Cell(4, 4).Range.Text = "41+51=92"
$t.Cell(4, 5).Range.Text = "8-7=1"
$t.Cell(5, 1).Range.Text = "21+65=86"
$t.Cell(5, 2).Range.Text = "61-44=17"
$t.Cell(5, 3).Range.Text = "68+1=69"
$t.Cell(5, 4).Range.Text = "21+71=92"
$t.Cell(5, 5).Range.Text = "79-69=10"
$t.Cell(6, 1).Range.Text = "27-15=12"
$t.Cell(6, 2).Range.Text = "19+33=52"
$t.Cell(6, 3).Range.Text = "22-21=1"
$t.Cell(6, 4).Range.Text = "27-7=20"
$t.Cell(6, 5).Range.Text = "53+20=73"
$t.Cell(7, 1).Range.Text = "25+60=85"
$t.Cell(7, 2).Range.Text = "19+63=82"
$t.Cell(7, 3).Range.Text = "85-0=85"
$t.Cell(7, 4).Range.Text = "86-4=82"
$t.Cell(7, 5).Range.Text = "16+74=90"
$t.Cell(8, 1).Range.Text = "10+41=51"
$t.Cell(8, 2).Range.Text = "29+44=73"
$t.Cell(8, 3).Range.Text = "64+24=88"
$t.Cell(8, 4).Range.Text = "62-61=1"
$t.Cell(8, 5).Range.Text = "35+23=58"
$t.Cell(9, 1).Range.Text = "86-45=41"
$t.Cell(9, 2).Range.Text = "34+34=68"
$t.Cell(9, 3).Range.Text = "36-4=32"
$t.Cell(9, 4).Range.Text = "7+20=27"
$t.Cell(9, 5).Range.Text = "95-84=11"
$t.Cell(10, 1).Range.Text = "73-21=52"
$t.Cell(10, 2).Range.Text = "79-14=65"
$t.Cell(10, 3).Range.Text = "77-69=8"
$t.Cell(10, 4).Range.Text = "38+15=53"
$t.Cell(10, 5).Range.Text = "1+27=28"
$t.Cell(11, 1).Range.Text = "35+57=92"
$t.Cell(11, 2).Range.Text = "51-41=10"
$t.Cell(11, 3).Range.Text = "90+1=91"
$t.Cell(11, 4).Range.Text = "69-60=9"
$t.Cell(11, 5).Range.Text = "29-18=11"
$t.Cell(12, 1).Range.Text = "74-13=61"
$t.Cell(12, 2).Range.Text = "69-59=10"
$t.Cell(12, 3).Range.Text = "35+18=53"
$t.Cell(12, 4).Range.Text = "90-14=76"
$t.Cell(12, 5).Range.Text = "88-48=40"
$t.Cell(13, 1).Range.Text = "29+30=59"
$t.Cell(13, 2).Range.Text = "76+4=80"
$t.Cell(13, 3).Range.Text = "49+42=91"
$t.Cell(13, 4).Range.Text = "74-37=37"
$t.Cell(13, 5).Range.Text = "91-56=35"
$t.Cell(14, 1).Range.Text = "30-16=14"
$t.Cell(14, 2).Range.Text = "60-23=37"
$t.Cell(14, 3).Range.Text = "20+19=39"
$t.Cell(14, 4).Range.Text = "88+9=97"
$t.Cell(14, 5).Range.Text = "96-67=29"
$t.Cell(15, 1).Range.Text = "77-0=77"
$t.Cell(15, 2).Range.Text = "17-14=3"
$t.Cell(15, 3).Range.Text = "96-95=1"
$t.Cell(15, 4).Range.Text = "39+14=53"
$t.Cell(15, 5).Range.Text = "55-27=28"
$t.Cell(16, 1).Range.Text = "80-48=32"
$t.Cell(16, 2).Range.Text = "55-19=36"
$t.Cell(16, 3).Range.Text = "23+42=65"
$t.Cell(16, 4).Range.Text = "49-35=14"
$t.Cell(16, 5).Range.Text = "28+59=87"
$t.Cell(17, 1).Range.Text = "92-60=32"
$t.Cell(17, 2).Range.Text = "98-65=33"
$t.Cell(17, 3).Range.Text = "30+54=84"
$t.Cell(17, 4).Range.Text = "5+73=78"
$t.Cell(17, 5).Range.Text = "39-9=30"
$t.Cell(18, 1).Range.Text = "89-36=53"
$t.Cell(18, 2).Range.Text = "94-53=41"
$t.Cell(18, 3).Range.Text = "75+8=83"
$t.Cell(18, 4).Range.Text = "79-30=49"
$t.Cell(18, 5).Range.Text = "37-31=6"
$t.Cell(19, 1).Range.Text = "19+7=26"
$t.Cell(19, 2).Range.Text = "11+35=46"
$t.Cell(19, 3).Range.Text = "51+29=80"
$t.Cell(19, 4).Range.Text = "44+42=86"
$t.Cell(19, 5).Range.Text = "84+12=96"
$t.Cell(20, 1).Range.Text = "38+21=59"
$t.Cell(20, 2).Range.Text = "84+7=91"
$t.Cell(20, 3).Range.Text = "15+17=32"
$t.Cell(20, 4).Range.Text = "31-6=25"
$t.Cell(20, 5).Range.Text = "38-30=8"

Write-Output "done"
